$wb = $excel.ActiveWorkbook

$wsJob = $wb.Worksheets.Item("Job to Run")
$wsAll = $wb.Worksheets.Item("All")

# Append the new job entry as a new row on the "All" sheet
$wsAll.Activate()
$wsAll.Range("A5").Value = "G1_MC_RCG_PNAS_pos_2021"
$wsAll.Range("B5").Value = 4
$wsAll.Range("C5").Value = 4
$wsAll.Range("D5").Value = "POS"

# Copy formatting (style) from row 2 to the new row 5
$wsAll.Range("A2:D2").Copy()
$wsAll.Range("A5:D5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Select B12 on "All" sheet to match the saved selection state
$wsAll.Range("B12").Select()

# Update "Job to Run" sheet row 2 with the new job entry
$wsJob.Activate()
$wsJob.Range("A2").Value = "G1_MC_RCG_PNAS_pos_2021"
$wsJob.Range("B2").Value = 4
$wsJob.Range("C2").Value = 4
$wsJob.Range("D2").Value = "POS"

# Select A2:D2 on "Job to Run" sheet to match the saved selection state
$wsJob.Range("A2:D2").Select()
